$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 ("Contact" / "No display for ContactDetail") becomes
# ("Jurisdiction" / "United States of America")
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# it is removed entirely, shifting all following rows up by one.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" updates ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition now reflect the profile's own
# title & description instead of the generic Extension text.
$elements.Range("K2").Value = "Measure Report Evidence Value"
$elements.Range("L2").Value = "Output value of rule definition within a measure report"
